# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# The existing sheet runs from A1:AB50 (a header row in row 1, plus 49 data
# rows). We extend it with three new trailing columns: AC (Wins), AD
# (Losses), AE (Ties). The header cells get the same bold/bordered/centered
# style already used for the other header cells (copied from AB1, the last
# existing header), and every data row (2-50) gets the season record
# 84-78-0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header formatting from the last existing header cell (AB1) onto
# the three new header cells so they match the rest of row 1 (bold, thin
# border, centered/top aligned).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# New header labels.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Season record values for every player row.
$ws.Range("AC2:AC50").Value = 84
$ws.Range("AD2:AD50").Value = 78
$ws.Range("AE2:AE50").Value = 0
